$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting so values
# like "1.00" or "67.233.19" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.233.19'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '3.836.40'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '449.42'
$ws.Range('E5').Value = '  +6.83%  '
$ws.Range('D6').Value = '146.84'
$ws.Range('E6').Value = '  +13.99%  '
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +3.79%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').Value = '  +3.21%  '
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('D11').Value = '0.0000325'
$ws.Range('E11').Value = '  -8.48%  '
$ws.Range('D12').Value = '43.71'
$ws.Range('E12').Value = '  +9.35%  '
$ws.Range('D13').Value = '10.37'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').Value = '4.454.76'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').Value = '14.94'
$ws.Range('E15').Value = '  -4.41%  '
$ws.Range('D16').Value = '3.866.41'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '20.10'
$ws.Range('E18').Value = '  +3.36%  '
$ws.Range('D19').Value = '1.16'
$ws.Range('E19').Value = '  +7.70%  '
$ws.Range('D20').Value = '67.347.78'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').Value = '426.10'
$ws.Range('E21').Value = '  +5.16%  '
$ws.Range('D22').Value = '14.78'
$ws.Range('E22').Value = '  +3.81%  '
$ws.Range('D23').Value = '3.24'
$ws.Range('E23').Value = '  +8.37%  '
$ws.Range('D24').Value = '86.47'
$ws.Range('E24').Value = '  +3.60%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '3.46'
$ws.Range('E25').Value = '  +8.71%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '37.38'
$ws.Range('E26').Value = '  +1.01%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +19.34%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').Value = '5.52'
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').Value = '9.82'
$ws.Range('E29').Value = '  +4.97%  '
$ws.Range('D30').Value = '750.89'
$ws.Range('E30').Value = '  +6.13%  '
$ws.Range('D31').Value = '13.78'
$ws.Range('E31').Value = '  +12.03%  '
$ws.Range('E32').Value = '  +11.81%  '
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('D34').Value = '43.56'
$ws.Range('E34').Value = '  +13.74%  '
$ws.Range('D35').Value = '0.157'
$ws.Range('E35').Value = '  +4.53%  '
$ws.Range('D36').Value = '57.08'
$ws.Range('E36').Value = '  +4.05%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +12.42%  '
$ws.Range('D39').Value = '0.0478'
$ws.Range('E39').Value = '  +5.75%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = '3.03'
$ws.Range('E40').Value = '  +2.68%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0715'
$ws.Range('E41').Value = '  -8.08%  '
$ws.Range('D42').Value = '0.344'
$ws.Range('E42').Value = '  +17.43%  '
$ws.Range('D43').Value = '0.140'
$ws.Range('E43').Value = '  +4.96%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '3.46'
$ws.Range('E45').Value = '  +4.37%  '
$ws.Range('D46').Value = '2.15'
$ws.Range('E46').Value = '  +6.07%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '2.48'
$ws.Range('E47').Value = '  +11.69%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.20'
$ws.Range('E48').Value = '  +3.45%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = '2.68'
$ws.Range('E49').Value = '  +5.24%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '145.64'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').Value = '2.90'
$ws.Range('E51').Value = '  +5.57%  '
